$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column Q (2020) mirrors the existing column P (2019) styling/format -
# duplicate column P's cell formatting into the new column Q.
$ws.Columns("P").Copy()
$ws.Columns("Q").Insert()

$ws.Range("Q4").Value = 2020

$ws.Range("Q5").Value = 99.3
$ws.Range("Q6").Value = 99.371420589467803
$ws.Range("Q7").Value = 99.319469393395053
$ws.Range("Q8").Value = 99.442213297634979
$ws.Range("Q9").Value = 98.766881972988841
$ws.Range("Q10").Value = 99.212798374809537
$ws.Range("Q11").Value = 99.799160124155549
$ws.Range("Q12").Value = 99.146991622239156
$ws.Range("Q13").Value = 99.538370126605429
$ws.Range("Q14").Value = 99.765563948945029

# Selection moves to P7 (matches the recorded sheetView selection in the edit).
$ws.Range("P7").Select()
